$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 15:16"

# Alemania (Germany) - updated case counts
$ws.Range("B8").Value = 31370
$ws.Range("C8").Value = 2314
$ws.Range("E8").Value = 30489

# Austria - updated case counts
$ws.Range("B15").Value = 4926
$ws.Range("C15").Value = 452
$ws.Range("E15").Value = 4892
$ws.Range("F15").Value = 20

# Portugal - updated case counts
$ws.Range("E18").Value = 2310
$ws.Range("F18").Value = 48
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 30

# Vietnam overtakes Jordania and Albania in the ranking (rows 82-84)
$ws.Range("A82").Value = "Vietnam"
$ws.Range("B82").Value = 132
$ws.Range("C82").Value = 9
$ws.Range("D82").Value = 17
$ws.Range("E82").Value = 115
$ws.Range("F82").Value = 3
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 0

$ws.Range("A83").Value = "Jordania"
$ws.Range("B83").Value = 127
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 1
$ws.Range("E83").Value = 126
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 0

$ws.Range("A84").Value = "Albania"
$ws.Range("B84").Value = 123
$ws.Range("C84").Value = 19
$ws.Range("D84").Value = 10
$ws.Range("E84").Value = 108
$ws.Range("F84").Value = 2
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 5

# Sri Lanka - updated case counts
$ws.Range("B92").Value = 102
$ws.Range("C92").Value = 5
$ws.Range("E92").Value = 100

# Camboya (Cambodia) - updated case counts
$ws.Range("B94").Value = 91
$ws.Range("C94").Value = 4
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = 87

# Bahamas / Suazilandia / Guinea / Curazao move above Groenlandia (rows 161-165);
# Groenlandia's own numbers are also revised
$ws.Range("A161").Value = "Bahamas"
$ws.Range("B161").Value = 4
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0

$ws.Range("A162").Value = "Suazilandia"
$ws.Range("B162").Value = 4
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 0
$ws.Range("E162").Value = 4
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("A163").Value = "Guinea"
$ws.Range("B163").Value = 4
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 0
$ws.Range("E163").Value = 4
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 0

$ws.Range("A164").Value = "Curazao"
$ws.Range("B164").Value = 4
$ws.Range("C164").Value = 0
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 3
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

$ws.Range("A165").Value = "Groenlandia"
$ws.Range("B165").Value = 4
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 2
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0
